# Commit before configuring .env
# Search and cart is incomplete till now.
#
# Two product option strings on the "Products" sheet are edited to carry
# their price add-on (RAM upsell copy), and the "Products" sheet becomes
# the active/selected tab again (with its selection parked on E7),
# while "Addresses" is left selected at F7 but is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsProducts  = $wb.Worksheets.Item("Products")
$wsAddresses = $wb.Worksheets.Item("Addresses")

# "8 GB" (row 2, Bose QuietComfort 45's RAM) becomes an upsell-priced option.
$wsProducts.Range("E2").Value = "4GB [+`$20.00]"

# "8 GB" (row 4, Microsoft Surface Pro 9's RAM) becomes a different upsell
# option - note this used to share the same shared string as E2 but the two
# now diverge into independent text.
$wsProducts.Range("E4").Value = "8GB [+`$60.00]"

# Addresses keeps its prior selection (F7) - just record it while it's still
# the active sheet so the cursor position is preserved there.
$wsAddresses.Activate()
$wsAddresses.Range("F7").Select()

# Make Products the active sheet again and restore its own selection to E7
# (it previously was scrolled/selected at D5 with a frozen topLeftCell of B1;
# re-activating + selecting resets the scroll position too).
$wsProducts.Activate()
$wsProducts.Range("E7").Select()
